# Add the "tt" prefix in front of the numeric IMDb id inside every
# trailer_url (column D) value, e.g.
#   https://www.imdb.com/title/0068646/videogallery
#   -> https://www.imdb.com/title/tt0068646/videogallery

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Locate the trailer_url column from the header row instead of hard-coding
# a column letter, and find the real extent of the data.
$usedRange = $ws.UsedRange
$firstRow = $usedRange.Row
$lastRow  = $firstRow + $usedRange.Rows.Count - 1
$firstCol = $usedRange.Column
$lastCol  = $firstCol + $usedRange.Columns.Count - 1

$urlCol = 0
for ($c = $firstCol; $c -le $lastCol; $c++) {
    $header = $ws.Cells.Item($firstRow, $c).Text
    if ($header -eq "trailer_url") {
        $urlCol = $c
    }
}
if ($urlCol -eq 0) { $urlCol = 4 }  # fallback: column D

$changed = 0
for ($r = $firstRow + 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $urlCol)
    $val = $cell.Text
    if ($val -match "^(https?://www\.imdb\.com/title/)(\d+)(/videogallery)$") {
        $cell.Value = $Matches[1] + "tt" + $Matches[2] + $Matches[3]
        $changed++
    }
}

Write-Host "Updated $changed trailer_url cell(s) with the tt-prefixed IMDb id."
